# Update the workbook to match the new release:
#  - bump Version / Status / Date / Contact metadata on the "Metadata" sheet
#  - swap the two "Mapping" columns (AK/AL) on the "Elements" sheet so the
#    business-mapping column now comes before the RIM-mapping column

$wb = $excel.ActiveWorkbook

# --- Metadata sheet -------------------------------------------------------
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B3").Value  = "0.4.0-snapshot-1"               # Version
$meta.Range("B6").Value  = "draft"                           # Status
$meta.Range("B8").Value  = "2024-05-23T12:16:26+00:00"       # Date
$meta.Range("B10").Value = "ANS (https://esante.gouv.fr)"    # Contact

# --- Elements sheet: swap Mapping columns (AK <-> AL) ----------------------
$elements = $wb.Worksheets.Item("Elements")

$lastRow = $elements.UsedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $akCell = $elements.Cells.Item($r, 37)
    $alCell = $elements.Cells.Item($r, 38)
    $akVal = $akCell.Value()
    $alVal = $alCell.Value()
    $akCell.Value = $alVal
    $alCell.Value = $akVal
}

# Swap the stored column widths too, so column 37 (now the long business
# mapping text, ~75 chars wide) is wide and column 38 (now RIM mapping,
# ~25 chars wide) is narrow again.
$elements.Columns.Item(37).ColumnWidth = 74.15
$elements.Columns.Item(38).ColumnWidth = 24.15
